# Generate Report for Handback
# Populates the "Latest Target File", "Latest Handback File" and
# "Latest Handback DateTime" columns for the zh-cn and de-de handback
# reports, links the newly-filled-in target file cell back to the source
# markdown file, and refreshes the overall status on the Overview sheet.

$wb = $excel.ActiveWorkbook

$sourceMdName = "6d63531b-c377-4854-bb56-6b0546cfa32e.md"
$sourceMdUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90a3c1c40685dab7066e7b2c01eccc4134bf4f91/e2e/6d63531b-c377-4854-bb56-6b0546cfa32e.md"

# ---------------------------------------------------------------
# Status: the handoff is complete and both locales are now back in
# sync with the source (en-US) content. This status string is shared
# across the Overview summary columns and each locale's Status column.
# ---------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------
# zh-cn handback report
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("I2").Value = $sourceMdName
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $sourceMdUrl, "", "", $sourceMdName) | Out-Null

$zhcn.Range("J2").Value = "6d63531b-c377-4854-bb56-6b0546cfa32e.fd6b7621bac04a3b6b2391af348f717d54342167.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-05 13:13:21"

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(9).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------
# de-de handback report
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("I2").Value = $sourceMdName
$dede.Hyperlinks.Add($dede.Range("I2"), $sourceMdUrl, "", "", $sourceMdName) | Out-Null

$dede.Range("J2").Value = "6d63531b-c377-4854-bb56-6b0546cfa32e.fd6b7621bac04a3b6b2391af348f717d54342167.de-de.xlf"
$dede.Range("K2").Value = "2016-09-05 13:13:28"

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(9).ColumnWidth = 39.166666666666664
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
